# Update odds data on Sheet1 (rows 3 and 4) to reflect refreshed FlashScore values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 ---
$ws.Range("G3").Value  = 2.2
$ws.Range("I3").Value  = 3.6
$ws.Range("L3").Value  = 4.33
$ws.Range("M3").Value  = 1.1
$ws.Range("N3").Value  = 7
$ws.Range("O3").Value  = 1.44
$ws.Range("P3").Value  = 2.63
$ws.Range("Q3").Value  = 2.4
$ws.Range("R3").Value  = 1.53
$ws.Range("AF3").Value = 67
$ws.Range("AH3").Value = 9
$ws.Range("AR3").Value = 67
$ws.Range("AV3").Value = 5.5

# --- Row 4 ---
$ws.Range("G4").Value  = 1.85
$ws.Range("H4").Value  = 3.5
$ws.Range("I4").Value  = 4.33
$ws.Range("J4").Value  = 2.5
$ws.Range("K4").Value  = 2
$ws.Range("N4").Value  = 8
$ws.Range("O4").Value  = 1.44
$ws.Range("P4").Value  = 2.63
$ws.Range("Q4").Value  = 2.35
$ws.Range("R4").Value  = 1.57
$ws.Range("X4").Value  = 7.5
$ws.Range("Z4").Value  = 15
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 15
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 10
$ws.Range("AT4").Value = 9
$ws.Range("AV4").Value = 6
$ws.Range("AY4").Value = 101
$ws.Range("BA4").Value = 351
